$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 35; existing rows 35..88 shift down to 36..89.
$ws.Rows("35:35").Insert()

# The row that used to be 35 is now row 36. Duplicate its contents/format into
# the freshly inserted (currently blank) row 35.
$ws.Rows("36:36").Copy()
$ws.Rows("35:35").PasteSpecial()
$excel.CutCopyMode = $false

# Finally, update the date on the newly created row 35 to the new record's date.
$ws.Range("D35").Value2 = 44757
